$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp title (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 00:34"

# --- Reorder countries: Colombia now ranks above Serbia ---
# Row 44 used to be Serbia, row 45 used to be Colombia.
# After the update, Colombia (with new, higher case numbers) takes row 44,
# and Serbia (unchanged numbers) drops to row 45.
$ws.Range("A44").Value = "Colombia"
$ws.Range("B44").Value = 10051
$ws.Range("C44").Value = 595
$ws.Range("D44").Value = 2424
$ws.Range("E44").Value = 7199
$ws.Range("F44").Value = 129
$ws.Range("G44").Value = 21
$ws.Range("H44").Value = 428

$ws.Range("A45").Value = "Serbia"
$ws.Range("B45").Value = 9943
$ws.Range("C45").Value = 95
$ws.Range("D45").Value = 2453
$ws.Range("E45").Value = 7281
$ws.Range("F45").Value = 45
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 209

# --- Reorder countries: Tanzania moves above Jordania / Gabon / Malta / Jamaica ---
# Row 120 used to be Jordania, now becomes Tanzania (new, higher numbers).
# Jordania, Gabon, Malta and Jamaica each shift down one row (121-124),
# keeping their previous statistics.
$ws.Range("A120").Value = "Tanzania"
$ws.Range("B120").Value = 509
$ws.Range("C120").Value = 29
$ws.Range("D120").Value = 183
$ws.Range("E120").Value = 305
$ws.Range("F120").Value = 7
$ws.Range("G120").Value = 5
$ws.Range("H120").Value = 21

$ws.Range("A121").Value = "Jordania"
$ws.Range("B121").Value = 508
$ws.Range("C121").Value = 14
$ws.Range("D121").Value = 385
$ws.Range("E121").Value = 114
$ws.Range("F121").Value = 5
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 9

$ws.Range("A122").Value = "Gabon"
$ws.Range("B122").Value = 504
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 110
$ws.Range("E122").Value = 386
$ws.Range("F122").Value = 1
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 8

$ws.Range("A123").Value = "Malta"
$ws.Range("B123").Value = 489
$ws.Range("C123").Value = 3
$ws.Range("D123").Value = 419
$ws.Range("E123").Value = 65
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 5

$ws.Range("A124").Value = "Jamaica"
$ws.Range("B124").Value = 488
$ws.Range("C124").Value = 10
$ws.Range("D124").Value = 58
$ws.Range("E124").Value = 421
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 9

# --- Straightforward statistic updates (no reordering) ---

# Estados Unidos
$ws.Range("B4").Value = 1318787
$ws.Range("C4").Value = 26164
$ws.Range("E4").Value = 1018276
$ws.Range("G4").Value = 1575
$ws.Range("H4").Value = 78503

# Alemania
$ws.Range("B10").Value = 170588
$ws.Range("C10").Value = 1158
$ws.Range("E10").Value = 21378

# Argentina
$ws.Range("B57").Value = 5611
$ws.Range("C57").Value = 240
$ws.Range("E57").Value = 3659
$ws.Range("G57").Value = 11
$ws.Range("H57").Value = 293

# Ghana
$ws.Range("F62").Value = 8

# Uruguay
$ws.Range("B112").Value = 694
$ws.Range("C112").Value = 10
$ws.Range("D112").Value = 506
$ws.Range("E112").Value = 170
$ws.Range("F112").Value = 7
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 18

# Estado de Palestina
$ws.Range("D129").Value = 228
$ws.Range("E129").Value = 145
